$d = $word.ActiveDocument

# Locate the paragraph that holds the "Ver no Jupiter ..." text and the one
# right before it (the blank paragraph separating it from the requirements
# line) as well as the paragraph with the "(c) 2020 ..." footer text, then
# remove that whole block of three paragraphs in one go (including their
# paragraph marks).
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        # The blank paragraph right before this one is also removed.
        $startPara = $d.Paragraphs.Item($i - 1)
    }

    if ($t -like "*Contact: luizeleno@usp.br*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}
